$d = $word.ActiveDocument

# Update the date paragraph at the top of the document
$d.Content.Find.Execute("2024-02-29 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-03-01 Friday", 2)

# New values for each of the 100 table cells, in row-major reading order
$newValues = @(
    "25+73=98",
    "10+11=21",
    "15+41=56",
    "17+66=83",
    "58-45=13",
    "21+76=97",
    "71-18=53",
    "16+26=42",
    "18+59=77",
    "2+8=10",
    "82-75=7",
    "94-62=32",
    "88-25=63",
    "1+22=23",
    "63-60=3",
    "84-25=59",
    "98-61=37",
    "58-37=21",
    "56+43=99",
    "83-63=20",
    "89+3=92",
    "93-88=5",
    "45-42=3",
    "75-64=11",
    "70+7=77",
    "49-22=27",
    "13+86=99",
    "62+34=96",
    "32+60=92",
    "60-59=1",
    "69+19=88",
    "31+20=51",
    "5+23=28",
    "98-92=6",
    "84-23=61",
    "89-74=15",
    "83-14=69",
    "55+43=98",
    "63-53=10",
    "63-34=29",
    "23+22=45",
    "72-52=20",
    "6+9=15",
    "94-88=6",
    "57+13=70",
    "66-59=7",
    "70-70=0",
    "31+42=73",
    "15+68=83",
    "50+43=93",
    "86-66=20",
    "84-4=80",
    "88-72=16",
    "22+62=84",
    "85-1=84",
    "94-80=14",
    "33+30=63",
    "14+75=89",
    "50+26=76",
    "96-53=43",
    "76+9=85",
    "34-23=11",
    "47+4=51",
    "77-30=47",
    "54-24=30",
    "25+63=88",
    "11+74=85",
    "61-58=3",
    "87-29=58",
    "14-0=14",
    "16+38=54",
    "12+4=16",
    "16+22=38",
    "72-13=59",
    "87+5=92",
    "85-73=12",
    "89-37=52",
    "11+69=80",
    "90-90=0",
    "29+4=33",
    "0+22=22",
    "88-82=6",
    "16+13=29",
    "30+20=50",
    "48+31=79",
    "1+30=31",
    "96-20=76",
    "93-77=16",
    "71-58=13",
    "36-25=11",
    "76-38=38",
    "31-27=4",
    "89-59=30",
    "45+50=95",
    "61-37=24",
    "63-27=36",
    "69-0=69",
    "33-25=8",
    "0+90=90",
    "20+70=90"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count

$i = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i++
    }
}

Write-Output "Updated $i cells"
